$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.245278716087341
$ws.Range("B1").Value = 2.213175773620605
$ws.Range("C1").Value = 6.083687305450439
$ws.Range("D1").Value = 1.315924167633057
$ws.Range("E1").Value = 1.327541470527649
